$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 3 (pushes existing rows 3-11 down to 4-12,
# carrying their values/styles along), to make room for the missing
# "9:00 - 10:00" interval.
$ws.Rows("3:3").Insert()

# Fill in the new interval label in column A of the newly inserted row.
$ws.Range("A3").Value2 = "9:00 - 10:00"

# Match the new selection left behind in the saved workbook.
[void]$ws.Range("A3").Select()
